$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 35 (doc-committee / Sondra Stegenga),
# which pushes it down to row 36 and makes room for a new "affiliation"
# entry (Data Visualization Society) at row 35.
$ws.Rows.Item(35).Insert()

# Populate the new row 35 with the Data Visualization Society affiliation.
$ws.Range("A35").Value = "affiliation"
$ws.Range("D35").Value = "Data Visualization Society"
$ws.Range("G35").Value = "https://www.datavisualizationsociety.com"

# Match the "what" column formatting (wrap text) and row height used by
# the other affiliation rows (D33, D34, D36).
$ws.Range("D35").WrapText = $true
$ws.Rows.Item(35).RowHeight = 17

# Update the view state to reflect where the author was working when the
# file was saved.
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("G35").Select()
